$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1746.6666
$ws.Range("I29").Value = 120
$ws.Range("K29").Value = 360
$ws.Range("M29").Value = -79

$ws.Range("H132").Value = 1892.5964
$ws.Range("I132").Value = 1347.5676
$ws.Range("J132").Value = 2900.9
$ws.Range("K132").Value = 4042.7028
$ws.Range("L132").Value = 8702.700000000001
$ws.Range("M132").Value = -1512.7028
$ws.Range("N132").Value = -13762.7

$ws.Range("H137").Value = 1330.4082
$ws.Range("I137").Value = 1208.5938
$ws.Range("J137").Value = 1559.7059
$ws.Range("K137").Value = 3625.7814
$ws.Range("L137").Value = 4679.1177
$ws.Range("M137").Value = -1075.7814
$ws.Range("N137").Value = -9779.117699999999

$ws.Range("H138").Value = 2873.676
$ws.Range("I138").Value = 1917.9333
$ws.Range("J138").Value = 4527.846
$ws.Range("K138").Value = 5753.7999
$ws.Range("L138").Value = 13583.538
$ws.Range("M138").Value = -613.7999
$ws.Range("N138").Value = -23863.538

$ws.Range("H141").Value = 5955.1333
$ws.Range("I141").Value = 2646.2593
$ws.Range("J141").Value = 35735
$ws.Range("K141").Value = 7938.777900000001
$ws.Range("L141").Value = 107205
$ws.Range("M141").Value = -2758.777900000001
$ws.Range("N141").Value = -117565


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 168373.17
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 168373.17
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 168373.17
$ws.Range("N2").Value = -168599.17
$ws.Range("M2").ClearContents()

$ws.Range("H17").Value = 2000
$ws.Range("I17").Value = 2000
$ws.Range("K17").Value = 2000
$ws.Range("M17").Value = -1827

$ws.Range("H32").Value = 9630.811
$ws.Range("I32").Value = 10582.221
$ws.Range("J32").Value = 5888.6
$ws.Range("K32").Value = 10582.221
$ws.Range("L32").Value = 5888.6
$ws.Range("M32").Value = -10295.221
$ws.Range("N32").Value = -6462.6

$ws.Range("H74").Value = 739.9787
$ws.Range("I74").Value = 689.2
$ws.Range("J74").Value = 1030.1428
$ws.Range("K74").Value = 689.2
$ws.Range("L74").Value = 1030.1428
$ws.Range("M74").Value = 184.8
$ws.Range("N74").Value = -2778.1428

$ws.Range("H77").Value = 739.9787
$ws.Range("I77").Value = 689.2
$ws.Range("J77").Value = 1030.1428
$ws.Range("K77").Value = 3446
$ws.Range("L77").Value = 5150.714
$ws.Range("M77").Value = 922
$ws.Range("N77").Value = -13886.714

$ws.Range("H114").Value = 29750
$ws.Range("J114").Value = 29750
$ws.Range("L114").Value = 29750
$ws.Range("N114").Value = -38428

$ws.Range("H116").Value = 168373.17
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 168373.17
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 168373.17
$ws.Range("N116").Value = -172961.17
$ws.Range("M116").ClearContents()

$ws.Range("H134").Value = 44209.5
$ws.Range("J134").Value = 44209.5
$ws.Range("L134").Value = 44209.5
$ws.Range("N134").Value = -54349.5


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 168373.17
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 168373.17
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 168373.17
$ws.Range("N3").Value = -168601.17
$ws.Range("M3").ClearContents()

$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 5000
$ws.Range("K16").Value = 5000
$ws.Range("M16").Value = -4830


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1738.8518
$ws.Range("I31").Value = 1226.1464
$ws.Range("J31").Value = 3355.8462
$ws.Range("K31").Value = 1226.1464
$ws.Range("L31").Value = 3355.8462
$ws.Range("M31").Value = -931.1464000000001
$ws.Range("N31").Value = -3945.8462

$ws.Range("H34").Value = 1738.8518
$ws.Range("I34").Value = 1226.1464
$ws.Range("J34").Value = 3355.8462
$ws.Range("K34").Value = 1226.1464
$ws.Range("L34").Value = 3355.8462
$ws.Range("M34").Value = -1024.1464
$ws.Range("N34").Value = -3759.8462

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H58").Value = 598359.3
$ws.Range("I58").Value = 975476.1
$ws.Range("K58").Value = 975476.1
$ws.Range("M58").Value = -975273.1

$ws.Range("H60").Value = 38900
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 38900
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 38900
$ws.Range("N60").Value = -39922
$ws.Range("M60").ClearContents()

$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H118").Value = 24500
$ws.Range("J118").Value = 24500
$ws.Range("L118").Value = 24500
$ws.Range("N118").Value = -27814

$ws.Range("H134").Value = 1434.4906
$ws.Range("I134").Value = 1076.3334
$ws.Range("K134").Value = 3229.0002
$ws.Range("M134").Value = -694.0001999999999

$ws.Range("H136").Value = 598359.3
$ws.Range("I136").Value = 975476.1
$ws.Range("K136").Value = 2926428.3
$ws.Range("M136").Value = -2923878.3


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 13007.5
$ws.Range("I124").Value = 3030
$ws.Range("J124").Value = 16333.333
$ws.Range("K124").Value = 9090
$ws.Range("L124").Value = 48999.999
$ws.Range("M124").Value = -4180
$ws.Range("N124").Value = -58819.999

$ws.Range("H131").Value = 18538806
$ws.Range("J131").Value = 22749610
$ws.Range("L131").Value = 68248830
$ws.Range("N131").Value = -68258910


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3574.9722
$ws.Range("I132").Value = 2854.8965
$ws.Range("J132").Value = 6558.143
$ws.Range("K132").Value = 8564.6895
$ws.Range("L132").Value = 19674.429
$ws.Range("M132").Value = -6034.6895
$ws.Range("N132").Value = -24734.429


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1100.64
$ws.Range("I136").Value = 848.5217
$ws.Range("K136").Value = 2545.5651
$ws.Range("M136").Value = 4.434900000000198

$ws.Range("H140").Value = 15500
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 15500
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 15500
$ws.Range("N140").Value = -25860
$ws.Range("M140").ClearContents()

